$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2")

$categories = @(
    "1 = COVID-19",
    "2 = Public health messaging/gov't handling of COVID",
    "3 = Healthcare (access to care, short supply)",
    "4 = Long term care",
    "5 = Chronic disease (cancer, heart disease)",
    "6 = Mental health",
    "7 = Access to housing and food",
    "8 = Drug abuse",
    "9 = Inequality",
    "10 = Economy",
    "11 = Corrupt gov't",
    "12 = Climate change/environmental",
    "13 = Abortion",
    "14 = Reliance on meat",
    "15 = Domestic abuse",
    "16 = Misinformation",
    "17 = Internet addiction",
    "18 = Don't know"
)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $categories[$i]
}

$ws.Range("I1").Value = "index"

$ws.Columns.Item(1).ColumnWidth = 51.21875
$ws.Columns.Item(9).ColumnWidth = 13.5546875

$ws.Range("G8").Select()
